# Add results for the "lattice 2step" prefit experiment.
#
# Table 1 "SIS (lattice,L=100)" (header row 2, data starts row 3): add two
# new trailing columns J ("mf two step") and K ("mf two step prefit (1e5 on
# lattice data)") with a value only for the A=0 data row (row 3).
#
# Table 4 "SIS (nn,L=100)" (header row 21): add the same new trailing
# header "mf two step prefit (1e5 on lattice data)" in column K (no data
# rows filled in below it).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters for shared-string de-dup ordering: enter the K3 value, then
# the new header text (K2), then the J3 value, then the (already-existing)
# J2 header text, then reuse the header text again for K21.
$ws.Range("K3").Value = "0.35(0.02)"
$ws.Range("K2").Value = "mf two step prefit (1e5 on lattice data)"
$ws.Range("J3").Value = "0.45(0.012)"
$ws.Range("J2").Value = "mf two step"
$ws.Range("K21").Value = "mf two step prefit (1e5 on lattice data)"

# Scroll back to the top of the sheet and select J4, matching the saved
# view state in the workbook.
$ws.Range("J4").Select() | Out-Null
